$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - RandomForestRegressor (name unchanged), update metric values
$ws.Range("B3").Value = 0.02924799644477622
$ws.Range("C3").Value = 0.0313402629370561
$ws.Range("D3").Value = 0.05462273480529092

# Row 4 - rename model and update metric values
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.03224366651710333
$ws.Range("C4").Value = 0.03429022045241648
$ws.Range("D4").Value = 0.08101672222199238

# Row 5 - rename model and update metric values
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.02504122812745007
$ws.Range("C5").Value = 0.0252378950374292
$ws.Range("D5").Value = 0.02766330195057139
